# Apply odds updates per commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.3
$ws.Range("G2").Value = 3.4
$ws.Range("H2").Value = 2.24
$ws.Range("I2").Value = 2.28
$ws.Range("V2").Value = 1.78
$ws.Range("W2").Value = 1.41
$ws.Range("AM2").Value = 80
$ws.Range("AN2").Value = 30

# Row 3
$ws.Range("F3").Value = 1.96
$ws.Range("G3").Value = 2.04
$ws.Range("I3").Value = 4.5
$ws.Range("J3").Value = 3.65
$ws.Range("L3").Value = 1.37
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 3.85
$ws.Range("Q3").Value = 1.71
$ws.Range("R3").Value = 1.38
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.74
$ws.Range("U3").Value = 2.12
$ws.Range("V3").Value = 1.29
$ws.Range("W3").Value = 1.96
$ws.Range("X3").Value = 19.5
$ws.Range("Y3").Value = 20
$ws.Range("Z3").Value = 40
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 10.5
$ws.Range("AD3").Value = 21
$ws.Range("AE3").Value = 65
$ws.Range("AF3").Value = 15.5
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 22
$ws.Range("AI3").Value = 70
$ws.Range("AJ3").Value = 28
$ws.Range("AK3").Value = 25
$ws.Range("AL3").Value = 42
$ws.Range("AN3").Value = 16
$ws.Range("AO3").Value = 65

# Row 4
$ws.Range("F4").Value = 1.41
$ws.Range("G4").Value = 1.49
$ws.Range("I4").Value = 9.800000000000001
$ws.Range("K4").Value = 6.2
$ws.Range("N4").Value = 6.2
$ws.Range("O4").Value = 1.15
$ws.Range("P4").Value = 2.78
$ws.Range("Q4").Value = 1.45
$ws.Range("R4").Value = 1.73
$ws.Range("S4").Value = 2.12
$ws.Range("T4").Value = 1.68
$ws.Range("W4").Value = 3
$ws.Range("X4").Value = 40
$ws.Range("AA4").Value = 240
$ws.Range("AB4").Value = 15.5
$ws.Range("AC4").Value = 16
$ws.Range("AF4").Value = 13.5
$ws.Range("AG4").Value = 11.5
$ws.Range("AJ4").Value = 14.5
$ws.Range("AO4").Value = 90

# Row 5
$ws.Range("F5").Value = 2.14
$ws.Range("G5").Value = 2.34
$ws.Range("I5").Value = 3.9
$ws.Range("J5").Value = 3.2
$ws.Range("N5").Value = 3.25
$ws.Range("O5").Value = 1.35
$ws.Range("P5").Value = 1.75
$ws.Range("V5").Value = 1.35
$ws.Range("W5").Value = 1.74

# Row 6
$ws.Range("F6").Value = 1.81
$ws.Range("G6").Value = 1.91
$ws.Range("H6").Value = 4.7
$ws.Range("I6").Value = 5.5
$ws.Range("J6").Value = 3.55
$ws.Range("K6").Value = 3.8
$ws.Range("L6").Value = 1.37
$ws.Range("N6").Value = 3.2
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 1.75
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.28
$ws.Range("S6").Value = 3.85
$ws.Range("T6").Value = 1.98
$ws.Range("U6").Value = 1.81
$ws.Range("W6").Value = 2.08
$ws.Range("X6").Value = 15.5
$ws.Range("Y6").Value = 18
$ws.Range("AF6").Value = 10.5
$ws.Range("AH6").Value = 980
$ws.Range("AJ6").Value = 25
$ws.Range("AK6").Value = 1000
$ws.Range("AN6").Value = 1000

# Row 7
$ws.Range("G7").Value = 1.91
$ws.Range("Q7").Value = 1.79
$ws.Range("AG7").Value = 10

# Row 8
$ws.Range("F8").Value = 2.36
$ws.Range("I8").Value = 3.55
$ws.Range("K8").Value = 3.5
$ws.Range("P8").Value = 1.75
$ws.Range("Q8").Value = 2.12

# Row 9
$ws.Range("F9").Value = 1.95
$ws.Range("G9").Value = 2.04
$ws.Range("H9").Value = 4.3
$ws.Range("I9").Value = 4.7
$ws.Range("J9").Value = 3.4
$ws.Range("K9").Value = 3.75
$ws.Range("N9").Value = 3.25
$ws.Range("Q9").Value = 2.06
$ws.Range("R9").Value = 1.29
$ws.Range("S9").Value = 3.8
$ws.Range("V9").Value = 1.27
$ws.Range("W9").Value = 1.96
